$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of chat data to append (Fecha / Mensajes), starting at row 10
$rows = @(
    @("14-06-2021 03:45", "p de mensaje 2"),
    @("14-06-2021 03:52", "p de mensaje 2"),
    @("14-06-2021 03:58", "Paso 02"),
    @("14-06-2021 04:01", "Paso 02"),
    @("14-06-2021 04:02", "Paso 02"),
    @("14-06-2021 04:03", "hola prueba desde poooosmmaaaan"),
    @("14-06-2021 04:04", "Prueba - Mandando mensaje y agregando a la base de datos"),
    @("14-06-2021 04:06", "Prueba - Mandando mensaje y agregando a la base de datos"),
    @("14-06-2021 04:07", "Prueba - Mandando mensaje y agregando a la base de datos 2"),
    @("14-06-2021 04:07", "Prueba - Mandando mensaje y agregando a la base de datos 3"),
    @("14-06-2021 04:08", "Hola Bienvenido\n\nEste es un mensaje de Prueba. Deberás enviarnos el *numero* de lo que estas buscando.\n*1*. Opción 1\n*2*. Opción 2\n*3*. Opción 3\n\nDeberás enviar el *número* de la opción solicitado."),
    @("14-06-2021 04:08", "Hola Bienvenido`nEste es un mensaje de Prueba. Deberás enviarnos el *numero* de lo que estas buscando`n*1*. Opción 1`n*2*. Opción 2`n*3*. Opción 3`nDeberás enviar el *número* de la opción solicitado.")
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
